$wb = $excel.ActiveWorkbook

# Add the Deviations sheet at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Deviations"

# Header row
$ws.Range("A1").Value = "time"
$ws.Range("B1").Value = "S1in"
$ws.Range("C1").Value = "S2in"
$ws.Range("D1").Value = "Cin"
$ws.Range("E1").Value = "XTin"
$ws.Range("F1").Value = "Nin"

# Data rows
$ws.Range("A2").Value = 20
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

$ws.Range("A3").Value = 100
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1.2
$ws.Range("F3").Value = 1

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

$influent = $wb.Worksheets.Item("Influent")
$influent.Range("A1:E1").Select() | Out-Null

$ws.Range("F5").Select() | Out-Null
$ws.Activate() | Out-Null
